# feat: load tags when displaying notes and vocabulary entries
#
# Adds 9 new vocabulary entries to the ENGLISH sheet (rows 96-104) and
# 1 new note to the NOTES sheet (row 26).

$wb = $excel.ActiveWorkbook
$english = $wb.Worksheets.Item("ENGLISH")
$notes = $wb.Worksheets.Item("NOTES")

# Columns on ENGLISH: A=Word, B=Definition, C=Synonyms, D=Antonyms,
# E=Correct answer count, F=Created at
$newWords = @(
    @("unenviable", "", "unpleasant;difficult;undesirable", "", 0, "2021-11-16 12:53:53.320239"),
    @("vigilant", "", "watchful", "", 0, "2021-11-16 12:54:22.449368"),
    @("amplify", "", "expand;louden", "", 0, "2021-11-16 12:54:56.705918"),
    @("notorious", "well known, typically for some bad quality or deed", "infamous", "", 0, "2021-11-16 12:55:49.850794"),
    @("epitomize", "", "summarize;embody", "", 0, "2021-11-16 12:56:33.27447"),
    @("rigorous", "extremely thorough and careful", "strict;meticulous", "", 0, "2021-11-16 12:57:47.829373"),
    @("weed out", "", "isolate", "", 0, "2021-11-16 12:59:23.223598"),
    @("debrief", "question (someone, typically a soldier or spy) about a completed mission", "question", "", 0, "2021-11-16 13:00:53.826091"),
    @("infraction", "a violation of a law or agreement", "infringement", "", 0, "2021-11-16 13:01:47.533568")
)

$startRow = 96
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $entry = $newWords[$i]

    $english.Cells.Item($row, 1).Value = $entry[0]
    if ($entry[1] -ne "") {
        $english.Cells.Item($row, 2).Value = $entry[1]
    }
    $english.Cells.Item($row, 3).Value = $entry[2]
    if ($entry[3] -ne "") {
        $english.Cells.Item($row, 4).Value = $entry[3]
    }
    $english.Cells.Item($row, 5).Value = $entry[4]
    $english.Cells.Item($row, 6).Value = $entry[5]
}

# New note appended to the NOTES sheet
$notes.Cells.Item(26, 1).Value = "Less but better"
